$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column P (rows 3-33, the data rows that already
# have a "P" column) over to the new column Q so the new cells pick up the
# same cell styles (number format, font, borders, etc.) as their P
# neighbours.
$ws.Range("P3:P33").Copy() | Out-Null
$ws.Range("Q3:Q33").PasteSpecial(-4122) | Out-Null

# Header row: 2020
$ws.Range("Q3").Value = 2020

# Data rows: one new "2020" data point per row, matching column P's row.
$ws.Range("Q4").Value  = 1.9148453093736542
$ws.Range("Q5").Value  = 1.7453236044300597
$ws.Range("Q6").Value  = 2.0818900906859255
$ws.Range("Q7").Value  = 1.658050942694075
$ws.Range("Q8").Value  = 1.4467487937731931
$ws.Range("Q9").Value  = 1.8774124750304142
$ws.Range("Q10").Value = 0.96024351775610284
$ws.Range("Q11").Value = 0.63595936855594293
$ws.Range("Q12").Value = 1.2888424905592288
$ws.Range("Q13").Value = 1.6032353288937073
$ws.Range("Q14").Value = 2.4146715443031859
$ws.Range("Q15").Value = 0.79837132250209564
$ws.Range("Q16").Value = 1.3751327862596732
$ws.Range("Q17").Value = 0.67516929870164943
$ws.Range("Q18").Value = 2.1012817818869509
$ws.Range("Q19").Value = 1.5943738893736428
$ws.Range("Q20").Value = 1.5765365498500856
$ws.Range("Q21").Value = 1.6126194804433236
$ws.Range("Q22").Value = 0.37150276583809166
$ws.Range("Q23").Value = 0
$ws.Range("Q24").Value = 0.75125835774923
$ws.Range("Q25").Value = 2.8942542850468351
$ws.Range("Q26").Value = 2.72898263527357
$ws.Range("Q27").Value = 3.0545792215303034
$ws.Range("Q28").Value = 3.9473869708034344
$ws.Range("Q29").Value = 3.6031203021816895
$ws.Range("Q30").Value = 4.2520923837938582
$ws.Range("Q31").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("Q33").Value = 0

# Move the active selection to T1, matching the post-edit saved view state.
$ws.Range("T1").Select() | Out-Null
